$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the stored precision of the existing last reading (row 12)
$ws.Range("A12").Value = 45862.91689484954

# Append the new automated reading as row 13
$ws.Range("A13").Value = 45862.9585492305
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value = 2025
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 19.33
$ws.Range("E13").Value = 74.44
$ws.Range("F13").Value = 87.89
$ws.Range("G13").Value = 13.54
$ws.Range("H13").Value = "ESE"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "23:00:18"
